$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple cell value updates (no style changes) ---
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 7
$ws.Range("K2").Value = 8

$ws.Range("I6").Value = 8
$ws.Range("K6").Value = 6

$ws.Range("J8").Value = 4
$ws.Range("K8").Value = 6

$ws.Range("I12").Value = 3
$ws.Range("J12").Value = 5
$ws.Range("K12").Value = 8

$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 1
$ws.Range("I18").Value = 8
$ws.Range("J18").Value = 7
$ws.Range("K18").Value = 6

# --- Cells whose fill colour (and value) changes: L column turns orange (FFC000),
#     H column turns red (FF0000). Creation order below matches the order new
#     cellXfs entries must appear in (75, 76, 77). ---

$ws.Range("L15").Interior.Color = 49407
$ws.Range("L15").Value = 6

$ws.Range("L16").Interior.Color = 49407

$ws.Range("L17").Interior.Color = 49407
$ws.Range("L17").Value = 8

$ws.Range("H15").Interior.Color = 255
$ws.Range("H15").Value = 8

$ws.Range("H16").Interior.Color = 255
$ws.Range("H16").Value = 7

$ws.Range("H17").Interior.Color = 255
$ws.Range("H17").Value = 6

# --- Update active selection ---
[void]$ws.Range("N8").Select()
